$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.558.14'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.14%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.301.04'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.29%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.67'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '183.59'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -4.15%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.296.51'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.23%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.570'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.85%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.177'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.573'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.85%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.75'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.88%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.58%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '637.55'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.74%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.830.08'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.46'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.91%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.715.95'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.93%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.84'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.299.17'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.58%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.99'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.01%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.67'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '101.09'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.99'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.99%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.34%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.40'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.92%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '30.99'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.35'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -4.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.52'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '589.13'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +4.55%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -9.72%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '10.90'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.852.25'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.92%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.62'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.33%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.10%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.39'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +4.53%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '32.19'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -6.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.13'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -6.36%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -5.69%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.58%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -4.64%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -6.93%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.127'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.24%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '129.98'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +5.62%  '
